$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: add a bold/bordered header label to A1 of a sheet that already has
# a header style (style index 1 in the original file) applied to B1:E1 (or
# B1 alone). We clone the look of the existing header by copy/paste-special
# of formats from the neighboring header cell, then set the value.
# ---------------------------------------------------------------------------
function Set-HeaderCell($ws, $cellRef, $text, $formatSourceRef) {
    $ws.Range($formatSourceRef).Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range($cellRef).Value = $text
}

# ---------------------------------------------------------------------------
# Sheets 1-4: "Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio (MWMed)",
# "Atendimento a Ponta(MW)", "Potencia Incremental - SIN(MW)"
# Each one gets a new A1 header "Fonte/Tecnologia", loses the header style
# previously applied to A2:A12, and fixes accented labels.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)

    Set-HeaderCell $ws "A1" "Fonte/Tecnologia" "B1"

    $ws.Range("A2").Style = "Normal"
    $ws.Range("A2").Value = "Hidro"

    $ws.Range("A3").Style = "Normal"
    $ws.Range("A3").Value = "Gás Natural"

    $ws.Range("A4").Style = "Normal"
    $ws.Range("A4").Value = "Carvão"

    $ws.Range("A5").Style = "Normal"
    $ws.Range("A5").Value = "Nuclear"

    $ws.Range("A6").Style = "Normal"
    $ws.Range("A6").Value = "Óleos Comb"

    $ws.Range("A7").Style = "Normal"
    $ws.Range("A7").Value = "Biomassa"

    $ws.Range("A8").Style = "Normal"
    $ws.Range("A8").Value = "Eólica"

    $ws.Range("A9").Style = "Normal"
    $ws.Range("A9").Value = "Solar"

    $ws.Range("A10").Style = "Normal"
    $ws.Range("A10").Value = "Outros"

    $ws.Range("A11").Style = "Normal"
    $ws.Range("A11").Value = "Pot. Compl."

    $ws.Range("A12").Style = "Normal"
    $ws.Range("A12").Value = "GD"
}

# ---------------------------------------------------------------------------
# Sheet 5: "Emissoes Totais (MtCO2eq)"
# New A1 header "Período", relabel A2/A3, drop the header style from them,
# and remove row 4 ("Teto") entirely.
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

Set-HeaderCell $ws5 "A1" "Período" "B1"

$ws5.Range("A2").Style = "Normal"
$ws5.Range("A2").Value = "P.Médio"

$ws5.Range("A3").Style = "Normal"
$ws5.Range("A3").Value = "P.Crítico"

$ws5.Rows.Item(4).Delete() | Out-Null

# ---------------------------------------------------------------------------
# Sheet 6: "Custo Total (bilhões de R$)"
# New A1 header "Tipo Expansão", B1 becomes "2015" instead of "Custo",
# relabel A2/A3 and update their cost values.
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

Set-HeaderCell $ws6 "A1" "Tipo Expansão" "B1"

# B1 becomes the text label "2015" (matching the year-header cells used on
# the other sheets) rather than "Custo". Copy a real "2015" header cell
# (format + the underlying text value) so it lands as text, not a number.
$wb.Worksheets.Item(1).Range("B1").Copy() | Out-Null
$ws6.Range("B1").PasteSpecial(-4104) | Out-Null   # xlPasteAll

$ws6.Range("A2").Style = "Normal"
$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 592

$ws6.Range("A3").Style = "Normal"
$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99
